$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2020" column (Q) to the right of the existing "2019" column (P),
# reusing the formatting of the corresponding P-column cell for each row.

# Row 4 (year header): Q4 = 2020, same look as P4 (O4/N4/.../I4 all share this style)
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

# Row 5 (Education): Q5 = 53.2, same look as P5
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 53.2

# Row 6 (Health): Q6 = 23.2, same look as P6
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 23.2

# Row 7 (Social protection): Q7 = 10, same look as P7 but with a "0.0" number format
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 10
$ws.Range("Q7").NumberFormat = "0.0"

# Row 8 (total, bottom border): Q8 = 20, same look as P8 but with a "0.0" number format
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 20
$ws.Range("Q8").NumberFormat = "0.0"

# Match the saved selection/active cell of the source workbook.
[void]$ws.Range("P9").Select()
